# Add new row (row 13) to the "Card18" sheet, per commit "Add new row to Card18".
# Also normalizes D2/E2 (previously empty placeholders) to "nan", matching the
# rest of the sheet's empty-value convention, as happened in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# D2/E2 were the only "empty" cells in their columns; every other row uses the
# literal text "nan" as the empty-value sentinel. Bring them in line.
$ws.Cells.Item(2, 4).Value = "nan"
$ws.Cells.Item(2, 5).Value = "nan"

# New row 13 data. Format the cells that get a value as Text first so the
# numeric-looking values ("18", "150", "300") are stored as text, matching
# the sheet's existing convention (every other data cell in the sheet is
# text, not numeric). Columns D-K have no data in the new row (same as a
# blank/"nan"-less cell elsewhere in the sheet), so they're left untouched.
$row = 13
$ws.Range("A" + $row + ":C" + $row).NumberFormat = "@"
$ws.Range("L" + $row + ":N" + $row).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "18"
$ws.Cells.Item($row, 2).Value = "150"
$ws.Cells.Item($row, 3).Value = "300"
$ws.Cells.Item($row, 12).Value = "13\8\2025"
$ws.Cells.Item($row, 13).Value = "Dfk belt  947*2.5*1.5 قطع سير"
$ws.Cells.Item($row, 14).Value = "تم تغير سير+تم تشحيم الماكينه بالكامل"
